# Update countries & provincias Spain
# Applies the data refresh captured in the commit:
#  - Swap "Australia" / "Dinamarca" labels (Dinamarca's case count overtook Australia's,
#    so it moved above it in the sorted list) together with their updated stats.
#  - Refresh case counts for several other countries (Marruecos, Kuwait,
#    Bosnia y Herzegovina, Libano, Albania, Vietnam, Isla de Man).
#  - Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 33 / 34: Dinamarca overtakes Australia -------------------------
# Row 33 now shows Dinamarca's refreshed numbers.
$ws.Range("A33").Value = "Dinamarca"
$ws.Range("B33").Value = 6496
$ws.Range("C33").Value = 178
$ws.Range("D33").Value = 2235
$ws.Range("E33").Value = 3976
$ws.Range("F33").Value = 100
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 285

# Row 34 now shows Australia, keeping its previous numbers.
$ws.Range("A34").Value = "Australia"
$ws.Range("B34").Value = 6400
$ws.Range("C34").Value = 41
$ws.Range("D34").Value = 3598
$ws.Range("E34").Value = 2741
$ws.Range("F34").Value = 80
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 61

# --- Row 59: Marruecos ----------------------------------------------------
$ws.Range("B59").Value = 1838
$ws.Range("C59").Value = 75
$ws.Range("D59").Value = 210
$ws.Range("E59").Value = 1502

# --- Row 68: Kuwait --------------------------------------------------------
$ws.Range("B68").Value = 1355
$ws.Range("C68").Value = 55
$ws.Range("D68").Value = 176
$ws.Range("E68").Value = 1176
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 3

# --- Row 72: Bosnia y Herzegovina ------------------------------------------
$ws.Range("B72").Value = 1080
$ws.Range("C72").Value = 43
$ws.Range("E72").Value = 823

# --- Row 90: Libano ----------------------------------------------------------
$ws.Range("B90").Value = 641
$ws.Range("C90").Value = 9
$ws.Range("E90").Value = 540
$ws.Range("F90").Value = 33
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 21

# --- Row 97: Albania -------------------------------------------------------
$ws.Range("B97").Value = 475
$ws.Range("C97").Value = 8
$ws.Range("D97").Value = 248
$ws.Range("E97").Value = 203
$ws.Range("F97").Value = 5
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 24

# --- Row 114: Vietnam --------------------------------------------------------
$ws.Range("D114").Value = 167
$ws.Range("E114").Value = 98

# --- Row 115: Isla de Man -----------------------------------------------------
$ws.Range("D115").Value = 141
$ws.Range("E115").Value = 99

# --- Timestamp caption -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 11:52"
